$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E16: flip sign from negative to positive
$ws.Range("E16").Value = 60473972.810000002

# E18: replace static value with SUM formula
$ws.Range("E18").Formula = "=SUM(E12:E17)"

# E19: fill in previously blank cell
$ws.Range("E19").Value = 384700000

# E21: fill in previously blank cell with SUM formula
$ws.Range("E21").Formula = "=SUM(E18:E20)"

# E22: fill in previously blank cell
$ws.Range("E22").Value = 20015625

# E26: replace shared string " " with numeric value
$ws.Range("E26").Value = 1018613404

$wb.Save()
